# Update TPM-derived values for the Sostdc1-Lrp4 LR-pairs sheet
# (commit: "update scripts wuth new tpm")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new values for columns E,F,G,H,M,N,O,P,Q,R,S,T
$data = @{
    2 = @{ E=2; F=0.6666666666666666; G=0.21011; H=0.6303299999999999;
           M=0.665744;              N=1.997232;
           O=0.07058985944777574;   P=0.07058985944777574;
           Q=0.13987947184;         R=1.25891524656;
           S=0.07058985944777574;   T=0.07058985944777574 }
    3 = @{ E=2; F=0.6666666666666666; G=0.21011; H=0.6303299999999999;
           O=0.1615251915556908;    P=0.1615251915556908;
           Q=0.3200751306266666;    R=2.880676175639999;
           S=0.1615251915556908;    T=0.1615251915556908 }
    4 = @{ E=2; F=0.6666666666666666; G=0.21011; H=0.6303299999999999;
           M=1.786190333333334;     N=5.358571;
           O=0.1893925060938975;    P=0.1893925060938975;
           Q=0.3752964509366667;    R=3.37766805843;
           S=0.1893925060938975;    T=0.1893925060938975 }
    5 = @{ E=2; F=0.6666666666666666; G=0.21011; H=0.6303299999999999;
           M=2.606414666666667;     N=7.819244;
           O=0.2763621526932594;    P=0.2763621526932594;
           Q=0.5476337856133333;    R=4.928704070519999;
           S=0.2763621526932594;    T=0.2763621526932594 }
    6 = @{ E=2; F=0.6666666666666666; G=0.21011; H=0.6303299999999999;
           M=1.713409666666666;     N=5.140229;
           O=0.1816754601565471;    P=0.1816754601565471;
           Q=0.3600045050633333;    R=3.240040545569999;
           S=0.1816754601565471;    T=0.1816754601565471 }
    7 = @{ E=2; F=0.6666666666666666; G=0.21011; H=0.6303299999999999;
           M=1.136028333333333;     N=3.408085;
           O=0.1204548300528295;    P=0.1204548300528295;
           Q=0.2386909131166667;    R=2.14821821805;
           S=0.1204548300528295;    T=0.1204548300528295 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
